$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the 2 new columns (CodigoAgente, NUM_GRUPO) before the old "TIPOPOLIZA" column ---
$ws.Columns("F:G").Insert()

# --- Insert 1 new column (Es0km) before the old "Marca" column (after the shift above) ---
$ws.Columns("Q:Q").Insert()

# --- Headers for the newly inserted columns ---
$ws.Range("F1").Value = "CodigoAgente"
$ws.Range("G1").Value = "NUM_GRUPO"
$ws.Range("Q1").Value = "Es0km"

# --- New trailing column header (SinAsistenciaMecanica) ---
$ws.Range("Z1").Value = "SinAsistenciaMecanica"

# --- Style the blank F/G cells on rows 2-4 (centered, no border) like row 5/6 will be ---
$blankRange = $ws.Range("F2:G4")
$blankRange.HorizontalAlignment = -4108
$blankRange.Borders.LineStyle = -4142

# --- Row 2 ---
$ws.Range("Q2").Value = "si"
$ws.Range("P2").Value = 2021
$ws.Range("Z2").Value = "Sí"

# --- Row 3 ---
$ws.Range("Q3").Value = "no"
$ws.Range("P3").Value = 2021
$ws.Range("Z3").Value = "No"

# --- Row 4 (no change to P4, Q4 stays blank, no Z4) ---

# --- Row 5 ---
$ws.Range("F5").Value = 2302
$ws.Range("G5").Value = "Mattioli"
$ws.Range("N5").Value = "03/05/2022"
$ws.Range("P5").Value = 2021
$ws.Range("Q5").Value = "si"
$ws.Range("V5").Value = "RPR007"
$ws.Range("W5").Value = "ABC12RPR007"
$ws.Range("X5").Value = "ZAZ123RPR007"
$ws.Range("Z5").Value = "Sí"

# --- Row 6 ---
$ws.Range("F6").Value = 2302
$ws.Range("G6").Value = "Mattioli"
$ws.Range("N6").Value = "03/05/2022"
$ws.Range("P6").Value = 2021
$ws.Range("Q6").Value = "no"
$ws.Range("V6").Value = "RPR008"
$ws.Range("W6").Value = "ABC12RPR008"
$ws.Range("X6").Value = "ZAZ123RPR008"
$ws.Range("Z6").Value = "No"

# --- Style the F/G cells on rows 5-6 to match rows 2-4 (centered, no border) ---
$ws.Range("F5:G6").HorizontalAlignment = -4108
$ws.Range("F5:G6").Borders.LineStyle = -4142

# --- Sheet view: drop the frozen/topLeft scroll position, move selection to F6:G6 ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F6:G6").Select()
